$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 39x10 array covering A2:J40 with updated labels (col A) and metric values (cols B-J)
$data = New-Object 'object[,]' 39,10

$data[0,0] = 'dangerous_pass_completion_ratio_under_low_pressure'
$data[0,1] = 8.565383007793118
$data[0,2] = 6.532750709376883
$data[0,3] = 31.11
$data[0,4] = 1.371748502529763
$data[0,5] = 1.652472431035928
$data[0,6] = 6.399655955296732
$data[0,7] = 3.915117119233172
$data[0,8] = 9.892657187692578
$data[0,9] = 9.195543905844644

$data[1,0] = 'count_completed_dangerous_passes_under_low_pressure_per_match'
$data[1,1] = 1.385890456071857
$data[1,2] = 1.06381422774367
$data[1,3] = 30.28
$data[1,4] = 0.2396916952898746
$data[1,5] = 0.2783849758969473
$data[1,6] = 1.006043399606137
$data[1,7] = 0.6350296500358649
$data[1,8] = 1.611020134260564
$data[1,9] = 1.519382362112671

$data[2,0] = 'count_dangerous_pass_attempts_under_low_pressure_per_match'
$data[2,1] = 2.410640597050704
$data[2,2] = 2.000435193492698
$data[2,3] = 20.51
$data[2,4] = 0.3431588645066361
$data[2,5] = 0.3774369966919066
$data[2,6] = 1.899284956914031
$data[2,7] = 1.410419293421
$data[2,8] = 2.813676299716636
$data[2,9] = 2.578409050715873

$data[3,0] = 'count_completed_passes_under_low_pressure_per_match'
$data[3,1] = 15.31797612483565
$data[3,2] = 13.03439753150214
$data[3,3] = 17.52
$data[3,4] = 2.55816329417295
$data[3,5] = 2.326907235871071
$data[3,6] = 11.86438170919754
$data[3,7] = 8.309592396820982
$data[3,8] = 18.6267553478586
$data[3,9] = 17.86666724365256

$data[4,0] = 'count_completed_dangerous_passes_under_high_pressure_per_match'
$data[4,1] = 5.064005266963316
$data[4,2] = 4.315425007309899
$data[4,3] = 17.35
$data[4,4] = 0.8476403585258567
$data[4,5] = 0.7560831502999676
$data[4,6] = 3.891934750285112
$data[4,7] = 2.697585580340602
$data[4,8] = 5.95875779720528
$data[4,9] = 5.434765673650782

$data[5,0] = 'count_ball_retentions_under_low_pressure_per_match'
$data[5,1] = 17.6000761708881
$data[5,2] = 15.16546310725304
$data[5,3] = 16.05
$data[5,4] = 2.694819682605172
$data[5,5] = 2.529511414493091
$data[5,6] = 14.09506098416305
$data[5,7] = 10.17061116057284
$data[5,8] = 21.41058995292123
$data[5,9] = 20.06471417986028

$data[6,0] = 'count_pass_attempts_under_low_pressure_per_match'
$data[6,1] = 18.52614447702483
$data[6,2] = 16.12575585749938
$data[6,3] = 14.89
$data[6,4] = 2.645292757139543
$data[6,5] = 2.400074948502687
$data[6,6] = 14.83393502398734
$data[6,7] = 11.99711956248855
$data[6,8] = 21.97759441061222
$data[6,9] = 21.36224793103679

$data[7,0] = 'count_dangerous_pass_attempts_under_high_pressure_per_match'
$data[7,1] = 9.708868999212649
$data[7,2] = 8.461886087333903
$data[7,3] = 14.74
$data[7,4] = 1.265124375996538
$data[7,5] = 1.215034694327539
$data[7,6] = 7.832796240424177
$data[7,7] = 5.874420561207423
$data[7,8] = 11.35585064438947
$data[7,9] = 10.15048953285467

$data[8,0] = 'count_low_pressures_received_per_match'
$data[8,1] = 21.73595382610306
$data[8,2] = 19.21712725753227
$data[8,3] = 13.11
$data[8,4] = 2.768115433347392
$data[8,5] = 2.587869700880357
$data[8,6] = 17.8605774925761
$data[8,7] = 15.0903415185088
$data[8,8] = 25.57737837236299
$data[8,9] = 24.4053109544749

$data[9,0] = 'count_dangerous_pass_attempts_under_medium_pressure_per_match'
$data[9,1] = 7.005592049041718
$data[9,2] = 6.23119532036499
$data[9,3] = 12.43
$data[9,4] = 0.6883559467227141
$data[9,5] = 0.9770005713853567
$data[9,6] = 5.936851822843828
$data[9,7] = 4.447640376312369
$data[9,8] = 7.808242880329368
$data[9,9] = 7.825262055526816

$data[10,0] = 'count_completed_passes_under_medium_pressure_per_match'
$data[10,1] = 70.82781662781403
$data[10,2] = 64.08204753818748
$data[10,3] = 10.53
$data[10,4] = 11.70603751087383
$data[10,5] = 9.150788641712808
$data[10,6] = 51.87327252371745
$data[10,7] = 47.75615730667351
$data[10,8] = 80.28543360661705
$data[10,9] = 82.24582440433775

$data[11,0] = 'count_completed_dangerous_passes_under_medium_pressure_per_match'
$data[11,1] = 3.531582012538288
$data[11,2] = 3.198585986536426
$data[11,3] = 10.41
$data[11,4] = 0.4777068396430045
$data[11,5] = 0.526490897142727
$data[11,6] = 2.992255121770995
$data[11,7] = 2.164787592259651
$data[11,8] = 4.066142021508265
$data[11,9] = 4.234255617531347

$data[12,0] = 'count_ball_retentions_under_medium_pressure_per_match'
$data[12,1] = 82.06453223751136
$data[12,2] = 74.85314098456587
$data[12,3] = 9.630000000000001
$data[12,4] = 12.60318022223563
$data[12,5] = 9.703934034261355
$data[12,6] = 61.22153456743462
$data[12,7] = 57.41912897938928
$data[12,8] = 92.1648809940915
$data[12,9] = 94.02726667018787

$data[13,0] = 'pass_completion_ratio_under_low_pressure'
$data[13,1] = 48.20684748018029
$data[13,2] = 44.34708230359838
$data[13,3] = 8.699999999999999
$data[13,4] = 4.410336955218041
$data[13,5] = 5.41551182118277
$data[13,6] = 41.17080940203395
$data[13,7] = 31.85565092840364
$data[13,8] = 52.15025335521475
$data[13,9] = 51.55774323594751

$data[14,0] = 'count_pass_attempts_under_medium_pressure_per_match'
$data[14,1] = 86.77578665338726
$data[14,2] = 80.17453381111775
$data[14,3] = 8.23
$data[14,4] = 11.8235906267499
$data[14,5] = 8.68487993797738
$data[14,6] = 67.44513834778533
$data[14,7] = 64.67815926678404
$data[14,8] = 96.63709701407491
$data[14,9] = 97.11523198142409

$data[15,0] = 'count_medium_pressures_received_per_match'
$data[15,1] = 103.1371696892224
$data[15,2] = 96.31081883995293
$data[15,3] = 7.09
$data[15,4] = 12.9900778914515
$data[15,5] = 9.404844930823185
$data[15,6] = 81.31058055075862
$data[15,7] = 80.0968719601268
$data[15,8] = 113.8601288456971
$data[15,9] = 115.240759483247

$data[16,0] = 'count_completed_passes_under_high_pressure_per_match'
$data[16,1] = 67.72483866786148
$data[16,2] = 63.26187430073603
$data[16,3] = 7.05
$data[16,4] = 7.59987005470007
$data[16,5] = 6.553433011659468
$data[16,6] = 56.52968207690235
$data[16,7] = 53.52864686206652
$data[16,8] = 76.02431420612432
$data[16,9] = 77.78047154425167

$data[17,0] = 'ball_retention_ratio_under_low_pressure'
$data[17,1] = 50.69211739554204
$data[17,2] = 47.41532395453828
$data[17,3] = 6.91
$data[17,4] = 4.071726453903429
$data[17,5] = 5.38928054846813
$data[17,6] = 44.36828160627763
$data[17,7] = 34.94103122000642
$data[17,8] = 55.31308364809014
$data[17,9] = 54.34400387761803

$data[18,0] = 'dangerous_pass_completion_ratio_under_high_pressure'
$data[18,1] = 20.44579215155953
$data[18,2] = 19.19621465662603
$data[18,3] = 6.51
$data[18,4] = 2.71981633972312
$data[18,5] = 2.591374925295072
$data[18,6] = 17.14139831677985
$data[18,7] = 14.2514256866458
$data[18,8] = 24.26918550227828
$data[18,9] = 22.59032836611046

$data[19,0] = 'dangerous_pass_completion_ratio_under_medium_pressure'
$data[19,1] = 17.09933544338354
$data[19,2] = 16.11574960881844
$data[19,3] = 6.1
$data[19,4] = 1.562111033886586
$data[19,5] = 1.793241943431018
$data[19,6] = 15.35393226282491
$data[19,7] = 12.11857449498212
$data[19,8] = 19.08287982014093
$data[19,9] = 19.40650061270864

$data[20,0] = 'count_ball_retentions_under_high_pressure_per_match'
$data[20,1] = 95.52800501858587
$data[20,2] = 90.11001493881376
$data[20,3] = 6.01
$data[20,4] = 9.593761884252615
$data[20,5] = 7.667538627233863
$data[20,6] = 82.19198514296028
$data[20,7] = 79.07083020752353
$data[20,8] = 106.9130877683976
$data[20,9] = 107.4145261516224

$data[21,0] = 'difficult_pass_completion_ratio_under_medium_pressure'
$data[21,1] = 27.41551230945307
$data[21,2] = 28.82947715268287
$data[21,3] = -4.9
$data[21,4] = 1.308154954285623
$data[21,5] = 1.582716530604583
$data[21,6] = 25.8586329953945
$data[21,7] = 25.67136560734314
$data[21,8] = 29.24262009233778
$data[21,9] = 30.50855872867531

$data[22,0] = 'count_pass_attempts_under_high_pressure_per_match'
$data[22,1] = 86.44445893810969
$data[22,2] = 82.58892288068476
$data[22,3] = 4.67
$data[22,4] = 6.839652965167274
$data[22,5] = 5.885820093249933
$data[22,6] = 76.13021056849205
$data[22,7] = 75.7589486335914
$data[22,8] = 94.21542623186333
$data[22,9] = 97.11915244577327

$data[23,0] = 'count_high_pressures_received_per_match'
$data[23,1] = 129.2296056586938
$data[23,2] = 124.5727367736607
$data[23,3] = 3.74
$data[23,4] = 9.417453240133325
$data[23,5] = 7.090524919515879
$data[23,6] = 115.721445699044
$data[23,7] = 114.7203877108655
$data[23,8] = 140.5440402814107
$data[23,9] = 141.6616858668597

$data[24,0] = 'count_completed_difficult_passes_under_medium_pressure_per_match'
$data[24,1] = 8.190240281646732
$data[24,2] = 8.449312368547409
$data[24,3] = -3.07
$data[24,4] = 0.579144379160078
$data[24,5] = 0.6534429990879151
$data[24,6] = 7.314206004593141
$data[24,7] = 6.913255624881678
$data[24,8] = 8.88520852954443
$data[24,9] = 9.270697872223591

$data[25,0] = 'difficult_pass_completion_ratio_under_high_pressure'
$data[25,1] = 29.60760713956032
$data[25,2] = 28.92710311937263
$data[25,3] = 2.35
$data[25,4] = 2.123053123217991
$data[25,5] = 1.783374684312475
$data[25,6] = 26.87328754377431
$data[25,7] = 26.87943058574159
$data[25,8] = 32.36513245430616
$data[25,9] = 32.42490839628574

$data[26,0] = 'count_forced_losses_under_high_pressure_per_match'
$data[26,1] = 33.70160064010796
$data[26,2] = 34.46272183484692
$data[26,3] = -2.21
$data[26,4] = 0.6521277532702746
$data[26,5] = 2.252794687552074
$data[26,6] = 32.81138503953775
$data[26,7] = 30.71203724013567
$data[26,8] = 34.60878306036576
$data[26,9] = 38.42727754238924

$data[27,0] = 'difficult_pass_completion_ratio_under_low_pressure'
$data[27,1] = 9.388398995518534
$data[27,2] = 9.188153106081417
$data[27,3] = 2.18
$data[27,4] = 1.130937246287605
$data[27,5] = 1.745602047215163
$data[27,6] = 7.906618336989101
$data[27,7] = 7.023128530731448
$data[27,8] = 10.50731671795062
$data[27,9] = 12.27675402487855

$data[28,0] = 'count_forced_losses_under_low_pressure_per_match'
$data[28,1] = 4.135877655214959
$data[28,2] = 4.051664150279227
$data[28,3] = 2.08
$data[28,4] = 0.3524198999202106
$data[28,5] = 0.4554164450721871
$data[28,6] = 3.76551650841305
$data[28,7] = 3.445235941680335
$data[28,8] = 4.591894578775816
$data[28,9] = 4.978570830286243

$data[29,0] = 'ball_retention_ratio_under_medium_pressure'
$data[29,1] = 66.33539009640845
$data[29,2] = 65.07082769935018
$data[29,3] = 1.94
$data[29,4] = 2.737829022009495
$data[29,5] = 2.770607139227629
$data[29,6] = 61.97341069876638
$data[29,7] = 58.99554369833887
$data[29,8] = 68.6619255890492
$data[29,9] = 68.34417862337921

$data[30,0] = 'count_forced_losses_under_medium_pressure_per_match'
$data[30,1] = 21.07263745171105
$data[30,2] = 21.45767785538708
$data[30,3] = -1.79
$data[30,4] = 0.8682174740824362
$data[30,5] = 1.056154776632053
$data[30,6] = 20.089045983324
$data[30,7] = 19.54765754965999
$data[30,8] = 22.1216798325188
$data[30,9] = 23.57308190408416

$data[31,0] = 'count_completed_difficult_passes_under_high_pressure_per_match'
$data[31,1] = 8.771641647253919
$data[31,2] = 8.6196950759535
$data[31,3] = 1.76
$data[31,4] = 0.4822935827077808
$data[31,5] = 0.7272518798499724
$data[31,6] = 8.297084507028288
$data[31,7] = 7.343120046737932
$data[31,8] = 9.518549566911286
$data[31,9] = 9.74609971088778

$data[32,0] = 'ball_retention_ratio_under_high_pressure'
$data[32,1] = 63.79663852730737
$data[32,2] = 62.74636097489237
$data[32,3] = 1.67
$data[32,4] = 3.106221437029308
$data[32,5] = 2.259124687653503
$data[32,6] = 58.72040143786406
$data[32,7] = 58.42119873135636
$data[32,8] = 66.29717625903702
$data[32,9] = 66.85710831866621

$data[33,0] = 'pass_completion_ratio_under_high_pressure'
$data[33,1] = 65.26093577501423
$data[33,2] = 64.25097355208202
$data[33,3] = 1.57
$data[33,4] = 4.146226643238476
$data[33,5] = 2.724910069035285
$data[33,6] = 58.37236531261879
$data[33,7] = 58.33117335911105
$data[33,8] = 68.67364819100777
$data[33,9] = 68.3629216180415

$data[34,0] = 'count_completed_difficult_passes_under_low_pressure_per_match'
$data[34,1] = 1.616317785566582
$data[34,2] = 1.597280398074337
$data[34,3] = 1.19
$data[34,4] = 0.1794954536251064
$data[34,5] = 0.3017557285122426
$data[34,6] = 1.386414378774036
$data[34,7] = 1.171010686964805
$data[34,8] = 1.800102473785913
$data[34,9] = 2.176649273793562

$data[35,0] = 'count_difficult_pass_attempts_under_medium_pressure_per_match'
$data[35,1] = 16.88538319807508
$data[35,2] = 17.07793803699742
$data[35,3] = -1.13
$data[35,4] = 0.7756577961598865
$data[35,5] = 1.079689568576342
$data[35,6] = 16.36988532008793
$data[35,7] = 14.34622741423436
$data[35,8] = 18.17544044914571
$data[35,9] = 18.10816188367631

$data[36,0] = 'count_difficult_pass_attempts_under_low_pressure_per_match'
$data[36,1] = 3.185130243983958
$data[36,2] = 3.213056349755119
$data[36,3] = -0.87
$data[36,4] = 0.1501018409291491
$data[36,5] = 0.419047254270259
$data[36,6] = 3.010210159807008
$data[36,7] = 2.41342584488841
$data[36,8] = 3.377466693799771
$data[36,9] = 4.05566687084303

$data[37,0] = 'count_difficult_pass_attempts_under_high_pressure_per_match'
$data[37,1] = 18.54999338279208
$data[37,2] = 18.68158928812503
$data[37,3] = -0.7
$data[37,4] = 0.8505335732176285
$data[37,5] = 1.82301944320041
$data[37,6] = 17.50032242899635
$data[37,7] = 15.20830733644102
$data[37,8] = 19.86248050314585
$data[37,9] = 21.55051489425237

$data[38,0] = 'pass_completion_ratio_under_medium_pressure'
$data[38,1] = 65.38722993464907
$data[38,2] = 65.21654156881303
$data[38,3] = 0.26
$data[38,4] = 3.171738950890961
$data[38,5] = 3.145396453232943
$data[38,6] = 60.16462282370047
$data[38,7] = 58.42356704351844
$data[38,8] = 68.39181857179648
$data[38,9] = 69.28342128159461

$ws.Range("A2:J40").Value2 = $data

Write-Host "Done updating A2:J40"